$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing content (keeps formatting) so the shared-strings table
# is rebuilt from only the strings actually referenced after this edit,
# matching the diff (drops the now-unused "ECs" shared string).
$ws.Cells.ClearContents()

$headers = @(
  "Sending cluster",
  "Ligand symbol",
  "Receptor symbol",
  "Target cluster",
  "Ligand-expressing cells",
  "Ligand detection rate",
  "Ligand average expression value",
  "Ligand total expression value",
  "Ligand derived specificity of average expression value",
  "Ligand derived specificity of total expression value",
  "Receptor-expressing cells",
  "Receptor detection rate",
  "Receptor average expression value",
  "Receptor total expression value",
  "Receptor derived specificity of average expression value",
  "Receptor derived specificity of total expression value",
  "Edge average expression weight",
  "Edge total expression weight",
  "Edge average expression derived specificity",
  "Edge total expression derived specificity"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf16"
$ws.Range("C2").Value = "Fgfr4"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.678104
$ws.Range("H2").Value = 5.034312
$ws.Range("I2").Value = 0.7635196712427992
$ws.Range("J2").Value = 0.7635196712427992
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1817723333333333
$ws.Range("N2").Value = 0.5453170000000001
$ws.Range("O2").Value = 0.008966262009224884
$ws.Range("P2").Value = 0.008966262009224884
$ws.Range("Q2").Value = 0.305032879656
$ws.Range("R2").Value = 2.745295916904
$ws.Range("S2").Value = 0.006845917421560183
$ws.Range("T2").Value = 0.006845917421560183

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf16"
$ws.Range("C3").Value = "Fgfr4"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.678104
$ws.Range("H3").Value = 5.034312
$ws.Range("I3").Value = 0.7635196712427992
$ws.Range("J3").Value = 0.7635196712427992
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.09115
$ws.Range("N3").Value = 60.27345
$ws.Range("O3").Value = 0.9910337379907751
$ws.Range("P3").Value = 0.9910337379907752
$ws.Range("Q3").Value = 33.7150391796
$ws.Range("R3").Value = 303.4353526164
$ws.Range("S3").Value = 0.756673753821239
$ws.Range("T3").Value = 0.756673753821239

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Fgf16"
$ws.Range("C4").Value = "Fgfr4"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.519749
$ws.Range("H4").Value = 1.559247
$ws.Range("I4").Value = 0.2364803287572008
$ws.Range("J4").Value = 0.2364803287572008
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1817723333333333
$ws.Range("N4").Value = 0.5453170000000001
$ws.Range("O4").Value = 0.008966262009224884
$ws.Range("P4").Value = 0.008966262009224884
$ws.Range("Q4").Value = 0.09447598847766667
$ws.Range("R4").Value = 0.8502838962990001
$ws.Range("S4").Value = 0.0021203445876647
$ws.Range("T4").Value = 0.0021203445876647

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf16"
$ws.Range("C5").Value = "Fgfr4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.519749
$ws.Range("H5").Value = 1.559247
$ws.Range("I5").Value = 0.2364803287572008
$ws.Range("J5").Value = 0.2364803287572008
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 20.09115
$ws.Range("N5").Value = 60.27345
$ws.Range("O5").Value = 0.9910337379907751
$ws.Range("P5").Value = 0.9910337379907752
$ws.Range("Q5").Value = 10.44235512135
$ws.Range("R5").Value = 93.98119609215
$ws.Range("S5").Value = 0.2343599841695361
$ws.Range("T5").Value = 0.2343599841695361

